$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cell values (cryptos list refresh, GitHub Actions run)
$updates = @(
    @('D2', '23.046.30'),
    @('E2', '  -3.48%  '),
    @('D3', '1.601.38'),
    @('E3', '  -2.90%  '),
    @('E4', '  +0.07%  '),
    @('D6', '301.09'),
    @('E6', '  -2.94%  '),
    @('D7', '0.3772'),
    @('E7', '  -2.82%  '),
    @('D8', '0.3634'),
    @('E8', '  -4.99%  '),
    @('D9', '49.96'),
    @('E9', '  -2.39%  '),
    @('E10', '  -5.55%  '),
    @('E11', '  +0.09%  '),
    @('D12', '0.08116'),
    @('E12', '  -3.64%  '),
    @('E13', '  -4.36%  '),
    @('D14', '6.585'),
    @('E14', '  -5.73%  '),
    @('D15', '0.00001255'),
    @('E15', '  -4.38%  '),
    @('E16', '  -7.89%  '),
    @('D17', '1.596.50'),
    @('E17', '  -4.25%  '),
    @('D18', '91.91'),
    @('E18', '  -2.11%  '),
    @('D19', '0.06878'),
    @('E19', '  -1.30%  '),
    @('D20', '18.23'),
    @('E20', '  -6.34%  '),
    @('D21', '6.562'),
    @('E21', '  -5.34%  '),
    @('B22', 'Dai'),
    @('C22', 'https://coinranking.com/coin/MoTuySvg7+dai-dai'),
    @('D22', '1.002'),
    @('E22', '  +0.17%  '),
    @('B23', 'Cosmos'),
    @('C23', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'),
    @('D23', '13.11'),
    @('E23', '  -3.64%  '),
    @('B24', 'WrappedBTC'),
    @('C24', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'),
    @('D24', '23.054.31'),
    @('E24', '  -3.41%  '),
    @('B25', 'Toncoin'),
    @('C25', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'),
    @('D25', '2.368'),
    @('E25', '  -3.13%  '),
    @('B26', 'LidoDAOToken'),
    @('C26', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'),
    @('D26', '2.778'),
    @('E26', '  -4.47%  '),
    @('B27', 'EthereumClassic'),
    @('C27', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'),
    @('D27', '21.08'),
    @('E27', '  -3.82%  '),
    @('B28', 'Monero'),
    @('C28', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'),
    @('D28', '150.42'),
    @('E28', '  -1.63%  '),
    @('B29', 'HuobiToken'),
    @('C29', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'),
    @('D29', '5.258'),
    @('E29', '  -2.31%  '),
    @('B30', 'BitcoinCash'),
    @('C30', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'),
    @('D30', '133.00'),
    @('E30', '  -2.93%  '),
    @('B31', 'WEMIXTOKEN'),
    @('C31', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'),
    @('D31', '2.331'),
    @('E31', '  -6.07%  '),
    @('B32', 'Filecoin'),
    @('C32', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'),
    @('D32', '6.863'),
    @('E32', '  -10.87%  '),
    @('B33', 'WrappedliquidstakedEther2.0'),
    @('C33', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'),
    @('D33', '1.777.07'),
    @('E33', '  -2.24%  '),
    @('B34', 'ImmutableX'),
    @('C34', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'),
    @('D34', '0.9616'),
    @('E34', '  -2.59%  '),
    @('B35', 'Hedera'),
    @('C35', 'https://coinranking.com/coin/jad286TjB+hedera-hbar'),
    @('D35', '0.07663'),
    @('E35', '  -5.71%  '),
    @('B36', 'FraxShare'),
    @('C36', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'),
    @('D36', '10.45'),
    @('E36', '  -0.98%  '),
    @('B37', 'InternetComputer(DFINITY)'),
    @('C37', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'),
    @('D37', '6.288'),
    @('E37', '  -5.19%  '),
    @('B38', 'VeChain'),
    @('C38', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'),
    @('D38', '0.02708'),
    @('E38', '  -6.65%  '),
    @('B39', 'Algorand'),
    @('C39', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'),
    @('D39', '0.2535'),
    @('E39', '  -5.05%  '),
    @('B40', 'Stellar'),
    @('C40', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'),
    @('D40', '0.08884'),
    @('E40', '  -2.22%  '),
    @('B41', 'TrustWalletToken'),
    @('C41', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'),
    @('D41', '1.365'),
    @('E41', '  -3.69%  '),
    @('B42', 'TheSandbox'),
    @('C42', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'),
    @('D42', '0.7042'),
    @('E42', '  -6.43%  '),
    @('B43', 'Aptos'),
    @('C43', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'),
    @('D43', '12.58'),
    @('E43', '  -6.09%  '),
    @('B44', 'EnergySwap'),
    @('C44', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'),
    @('D44', '15.23'),
    @('E44', '  -8.40%  '),
    @('B45', 'Decentraland'),
    @('C45', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'),
    @('D45', '0.6614'),
    @('E45', '  -4.30%  '),
    @('B46', 'NEARProtocol'),
    @('C46', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'),
    @('D46', '2.319'),
    @('E46', '  -4.72%  '),
    @('B47', 'Frax'),
    @('C47', 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'),
    @('D47', '1.001'),
    @('E47', '  +0.08%  '),
    @('B48', 'PancakeSwap'),
    @('C48', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'),
    @('D48', '3.993'),
    @('E48', '  -2.38%  '),
    @('B49', 'Quant'),
    @('C49', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'),
    @('D49', '132.52'),
    @('E49', '  -0.83%  '),
    @('B50', 'Flow'),
    @('C50', 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'),
    @('D50', '1.238'),
    @('E50', '  +1.32%  '),
    @('B51', 'Cronos'),
    @('C51', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'),
    @('D51', '0.07920'),
    @('E51', '  -4.12%  ')
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $val = $u[1]
    $rng = $ws.Range($cellRef)
    # Force text storage so numeric-looking strings (prices, percentages)
    # are not auto-converted into numbers/dates by Excel.
    $rng.NumberFormat = "@"
    $rng.Value = $val
}
